$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data (but keep the header row) for columns F:J, rows 2-11
$ws.Range("F2:J11").ClearContents()
